$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E/J swap: "OrderStatus" (E1) is replaced by "Site" (was J1) ---
$ws.Range("E1").Value = "Site"
$ws.Range("J1").Clear()

# --- Row 2 clean-up: drop the stray record-id cells that accompanied the
#     removed OrderStatus/Site columns (A2:D2 held ids, E2 held "Closed",
#     J2 held another id). D2 keeps its distinctive style, so only its
#     contents are cleared; the rest are fully cleared (default style). ---
$ws.Range("A2").Clear()
$ws.Range("B2").Clear()
$ws.Range("C2").Clear()
$ws.Range("D2").ClearContents()
$ws.Range("E2").Clear()
$ws.Range("J2").Clear()

# --- Row 3 clean-up: same pattern, plus the old J3 id cell ---
$ws.Range("A3").Clear()
$ws.Range("B3").Clear()
$ws.Range("C3").Clear()
$ws.Range("D3").Clear()
$ws.Range("J3").Clear()

# --- Row 4 only held a stray id in A4; clearing it drops the whole row ---
$ws.Range("A4").Clear()

# --- New column Z: "workinghours" header + id value (plain/default style) ---
$ws.Range("Z1").Value = "workinghours"
$ws.Range("Z1").Style = "Normal"
$ws.Range("Z2").Value = "01mo0000000K7kC"
$ws.Range("Z2").Style = "Normal"

# --- Column width adjustments ---
$ws.Range("E1").EntireColumn.ColumnWidth = 20.83203125
$ws.Range("Z1").EntireColumn.ColumnWidth = 41.33203125

# --- View state: scrolled/selected around the new last column ---
$ws.Application.ActiveWindow.ScrollColumn = 22
$ws.Range("Z1:Z2").Select()
